$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("C2").Value = 47
$ws.Range("D2").Value = 0.0001680972054600716
$ws.Range("E2").Value = 0.05758977402001619
$ws.Range("F2").Value = 47
$ws.Range("G2").Value = 0.003439337480813265
$ws.Range("H2").Value = 0.005828536115586758
$ws.Range("I2").Value = 0.02147068968042731
$ws.Range("J2").Value = 0.02305149985477328
$ws.Range("K2").Value = 0.00104889739304781

$ws.Range("C3").Value = 47
$ws.Range("D3").Value = 0.002187512349337339
$ws.Range("E3").Value = 0.04136045090854168
$ws.Range("F3").Value = 47
$ws.Range("G3").Value = 0.002034189645200968
$ws.Range("H3").Value = 0.005864783655852079
$ws.Range("I3").Value = 0.01335992338135839
$ws.Range("J3").Value = 0.01779578160494566
$ws.Range("K3").Value = 0.0006804326549172401

$ws.Range("C4").Value = 47
$ws.Range("D4").Value = 0.003939507994800806
$ws.Range("E4").Value = 0.06558519415557384
$ws.Range("F4").Value = 47
$ws.Range("G4").Value = 0.003155517857521772
$ws.Range("H4").Value = 0.009292104747146368
$ws.Range("I4").Value = 0.02248858381062746
$ws.Range("J4").Value = 0.02733874786645174
$ws.Range("K4").Value = 0.0008859853260219097

$ws.Range("C5").Value = 47
$ws.Range("D5").Value = 0.0001551201567053795
$ws.Range("E5").Value = 0.03734809719026089
$ws.Range("F5").Value = 47
$ws.Range("G5").Value = 0.00212458148598671
$ws.Range("H5").Value = 0.003786101005971432
$ws.Range("I5").Value = 0.01385019673034549
$ws.Range("J5").Value = 0.01518407743424177
$ws.Range("K5").Value = 0.0006839144043624401

$ws.Range("D6").Value = 0.003559940028935671
$ws.Range("E6").Value = 0.1465741978026927
$ws.Range("G6").Value = 0.004158143885433674
$ws.Range("H6").Value = 0.01211459096521139
$ws.Range("I6").Value = 0.1046045836992562
$ws.Range("J6").Value = 0.02039917698130012
$ws.Range("K6").Value = 0.00143083930015564

$ws.Range("C8").Value = 47
$ws.Range("D8").Value = 0.0001680972054600716
$ws.Range("E8").Value = 0.05758977402001619
$ws.Range("F8").Value = 47
$ws.Range("G8").Value = 0.003439337480813265
$ws.Range("H8").Value = 0.005828536115586758
$ws.Range("I8").Value = 0.02147068968042731
$ws.Range("J8").Value = 0.02305149985477328
$ws.Range("K8").Value = 0.00104889739304781

$ws.Range("C9").Value = 47
$ws.Range("D9").Value = 0.002187512349337339
$ws.Range("E9").Value = 0.04136045090854168
$ws.Range("F9").Value = 47
$ws.Range("G9").Value = 0.002034189645200968
$ws.Range("H9").Value = 0.005864783655852079
$ws.Range("I9").Value = 0.01335992338135839
$ws.Range("J9").Value = 0.01779578160494566
$ws.Range("K9").Value = 0.0006804326549172401

$ws.Range("C10").Value = 47
$ws.Range("D10").Value = 0.003939507994800806
$ws.Range("E10").Value = 0.06558519415557384
$ws.Range("F10").Value = 47
$ws.Range("G10").Value = 0.003155517857521772
$ws.Range("H10").Value = 0.009292104747146368
$ws.Range("I10").Value = 0.02248858381062746
$ws.Range("J10").Value = 0.02733874786645174
$ws.Range("K10").Value = 0.0008859853260219097

$ws.Range("C11").Value = 47
$ws.Range("D11").Value = 0.0001551201567053795
$ws.Range("E11").Value = 0.03734809719026089
$ws.Range("F11").Value = 47
$ws.Range("G11").Value = 0.00212458148598671
$ws.Range("H11").Value = 0.003786101005971432
$ws.Range("I11").Value = 0.01385019673034549
$ws.Range("J11").Value = 0.01518407743424177
$ws.Range("K11").Value = 0.0006839144043624401

$ws.Range("D12").Value = 0.003559940028935671
$ws.Range("E12").Value = 0.1465741978026927
$ws.Range("G12").Value = 0.004158143885433674
$ws.Range("H12").Value = 0.01211459096521139
$ws.Range("I12").Value = 0.1046045836992562
$ws.Range("J12").Value = 0.02039917698130012
$ws.Range("K12").Value = 0.00143083930015564

